# Branchement t1b2 sur nouvelle BDD
# Rename header columns (row 1) to match the new database field names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "c1_1_8deg_final"
$ws.Range("E1").Value = "c1_2deg_final"
$ws.Range("F1").Value = "C1_final_date"
$ws.Range("H1").Value = "C1_initial_date"
$ws.Range("I1").Value = "c2_1_5deg_final"
$ws.Range("J1").Value = "c2_1_8deg_final"
$ws.Range("K1").Value = "c2_2deg_final"
$ws.Range("L1").Value = "C2_final_date"
$ws.Range("Q1").Value = "C2_initial_date"
